# edit.ps1 - applies the "medal totals table code matches" edit:
#  - Slide 1 ("G01-A" textbox): drop the stray trailing endParaRPr run
#  - Slide 1 (authors textbox): drop the stray trailing endParaRPr run
#  - Slide 12 ("Dataset type" / ":" ): merge the two runs into one run "Dataset type:"
#  - Slide 3 ("All Winners..." paragraph): reword/re-split into 4 runs
#  - Slide 3 ("Data " / "sample"): merge the two runs into one run "Data sample"
#  - Slide 3 (medal totals table): nudge the table right (x offset 395536 -> 611560)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1: title/author block
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)

# Shape 2 = "Text Placeholder 4" containing "G01-A" (single paragraph).
# The whole text frame only has this one paragraph, so deleting the entire
# TextRange and retyping the same text drops the stale trailing
# <a:endParaRPr> while keeping the run's own formatting (inherited from the
# deleted run) intact.
$gShape = $slide1.Shapes.Item(2)
$gTr = $gShape.TextFrame.TextRange
$gText = $gTr.Text
$gTr.Delete()
$gTr.Text = $gText

# Shape 3 = "Text Placeholder 4" containing the three author lines. Same
# trick, but re-typed with explicit paragraph breaks so the three
# paragraphs come back exactly as they were (minus the stray endParaRPr on
# the last one).
$aShape = $slide1.Shapes.Item(3)
$aTr = $aShape.TextFrame.TextRange
$aTr.Delete()
$aTr.Text = "70493 – Tiago Nascimento`r76102 – Miguel Cruz`r76394 – Daniel Trindade"

# ---------------------------------------------------------------------------
# Slide 12: "Dataset type" + ":" -> single run "Dataset type:"
# ---------------------------------------------------------------------------
$slide12 = $p.Slides.Item(12)
$s12Shape = $slide12.Shapes.Item(2)
$s12Tr = $s12Shape.TextFrame.TextRange
$s12Found = $s12Tr.Find("Dataset type")
$s12Merged = $s12Tr.Characters($s12Found.Start, $s12Found.Length + 1)
$s12Merged.Text = "Dataset type:"

# ---------------------------------------------------------------------------
# Slide 3: "Initial Dataset" content placeholder
# ---------------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$s3Shape = $slide3.Shapes.Item(2)
$s3Tr = $s3Shape.TextFrame.TextRange

# "All Winners – A " / "table containing all the podium finishes " / "of all time."
#   -> "All Winners – " / "contains " / "the podium finishes " / "of all time."
# Insert a placeholder run right after the first run, retarget it via Find
# (InsertAfter's own return value tracks the *source* range, not the new
# text, in this host) so it becomes a distinct run, then fix up the two
# run texts that actually changed.
$s3FirstRun = $s3Tr.Find("All Winners – A ")
$null = $s3FirstRun.InsertAfter("##TMP##")
$s3NewRun = $s3Tr.Find("##TMP##")
$s3NewRun.Text = "contains "
$s3FirstRunAgain = $s3Tr.Find("All Winners – A ")
$s3FirstRunAgain.Text = "All Winners – "
$s3ThirdRun = $s3Tr.Find("table containing all the podium finishes ")
$s3ThirdRun.Text = "the podium finishes "

# "Data " / "sample" -> single run "Data sample"
$s3DataFound = $s3Tr.Find("Data ")
$s3DataMerged = $s3Tr.Characters($s3DataFound.Start, $s3DataFound.Length + 6)
$s3DataMerged.Text = "Data sample"

# Medal totals table: shift right (x offset 395536 -> 611560 EMU; 1 pt = 12700 EMU)
$s3Table = $slide3.Shapes.Item(3)
$s3Table.Left = 611560 / 12700.0
